# ModbusMod.xlsx edit:
#   - deleted clap detection (SOUND_CLAPCOUNT / CLCNT* bitfield)
#   - changed hardware to version 01 (round pcb):
#       row 3  (STATUS):            fill previously-empty C3:K3 with "X" placeholders
#       row 23 (was SOUND_CLAPCOUNT): renamed to FAST_PIRCOUNT_SOUNDPERCENT and its
#                                     bitfield layout replaced with PIRCOUNT7..0 / X / SOUNDPERC6..0
#   - selection moved to B23

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 3: STATUS row gains "X" filler bits in columns C..K (previously blank) ---
$ws.Range("C3:K3").Value = "X"

# --- Row 23: rename the row label and rewrite its bit layout ---
$ws.Cells.Item(23, 2).Value = "FAST_PIRCOUNT_SOUNDPERCENT"   # B23

$row23 = @{
    3  = "PIRCOUNT7"   # C23
    4  = "PIRCOUNT6"   # D23
    5  = "PIRCOUNT5"   # E23
    6  = "PIRCOUNT4"   # F23
    7  = "PIRCOUNT3"   # G23
    8  = "PIRCOUNT2"   # H23
    9  = "PIRCOUNT1"   # I23
    10 = "PIRCOUNT0"   # J23
    11 = "X"           # K23
    12 = "SOUNDPERC6"  # L23
    13 = "SOUNDPERC5"  # M23
    14 = "SOUNDPERC4"  # N23
    15 = "SOUNDPERC3"  # O23
    16 = "SOUNDPERC2"  # P23
    17 = "SOUNDPERC1"  # Q23
    18 = "SOUNDPERC0"  # R23
}

foreach ($col in $row23.Keys) {
    $ws.Cells.Item(23, $col).Value = $row23[$col]
}

# --- Move / record the active selection on Sheet1 ---
$ws.Activate()
$ws.Range("B23").Select()
